$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 714, shifting existing rows 714:834 down to 715:835
$ws.Rows(714).Insert()

# Populate the newly inserted row 714 with the new record
$ws.Cells.Item(714, 1).Value = 6
$ws.Cells.Item(714, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(714, 3).Value = "Metropolitana"
$ws.Cells.Item(714, 4).Value = 45218
$ws.Cells.Item(714, 5).Value = 13
$ws.Cells.Item(714, 6).Value = 100112044
$ws.Cells.Item(714, 7).Value = "Perejil"
$ws.Cells.Item(714, 8).Value = "Sin especificar"
$ws.Cells.Item(714, 9).Value = "Primera"
$ws.Cells.Item(714, 10).Value = 130
$ws.Cells.Item(714, 11).Value = 17000
$ws.Cells.Item(714, 12).Value = 18000
$ws.Cells.Item(714, 13).Value = 17615
$ws.Cells.Item(714, 14).Value = "`$/docena de atados"
$ws.Cells.Item(714, 15).Value = "Región Metropolitana"
$ws.Cells.Item(714, 16).Value = 5872
$ws.Cells.Item(714, 17).Value = 3
$ws.Cells.Item(714, 18).Value = "Hortaliza"
